# Faculty view exams completed
# - Add a second completed-exam row: A3 = 2, B3 = email address (hyperlinked like B2)
# - Update the worksheet selection from A3 to C3

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data: serial number and faculty email address
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "17it051@charusat.edu.in"

# Turn the email into a mailto hyperlink, matching the existing B2 hyperlink
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:17it051@charusat.edu.in")

# Re-apply the original (Hyperlink) cell style so it matches B2's formatting
# exactly instead of the extra style the Add() call introduces
$ws.Range("B3").Style = $ws.Range("B2").Style

# Move the active selection to C3
$ws.Range("C3").Select()
